# Purchased another (bigger, cheaper) cooler - add a new line item to
# Juan's purchases sheet for a Peltier cooler, and push the Total row down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Juan's purchases")

# Insert a new row right above the existing "Total" row (row 10) - this
# shifts "Total" (and its SUM formula) down to row 11 automatically.
$ws.Rows.Item(10).Insert()

# New purchase: Component, Part #, Purpose, Buyer, Price
$ws.Range("A10").Value = "Peltier Cooler"
$ws.Range("B10").Value = "TEC1-12706"
$ws.Range("C10").Value = "Cooling"
$ws.Range("D10").Value = "Yohan"
$ws.Range("E10").Value = 10.48

# Match the currency formatting used by the other Price cells.
$ws.Range("E10").NumberFormat = $ws.Range("E9").NumberFormat

# Part # for the new component is italicized; the component name itself
# stays in the regular (non-bold) font used by the rest of the rows.
$ws.Range("A10").Font.Bold = $false
$ws.Range("B10").Font.Italic = $true

# The "Total" label (now on row 11) stays bold.
$ws.Range("A11").Font.Bold = $true

# Extend the SUM to cover the newly inserted row.
$ws.Range("E11").Formula = "=SUM(E2:E10)"

# Leave the selection where the author left it after the edit.
$null = $ws.Range("D14").Select()
